# Update nomination summary worksheet: relabel category rows with
# "<Category>, <Metric>" style labels, fix misspellings (Air Fotce ->
# Air Force, Marine Cotps -> Marine Corps), rename the second Civilian
# block to "Civilian (FS, PHS, CG, NOAA)", and rework the Summary block
# into four clearly named totals (new/carryover/confirmed/unconfirmed/
# withdrawn), dropping the old duplicate "Total withdrawn" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Labels'
$ws.Range("B1").Value = 'Values'
$ws.Range("A2").Value = 'Congress'
$ws.Range("B2").Value = 104
$ws.Range("A3").Value = 'Session'
$ws.Range("B3").Value = 2
$ws.Range("A6").Value = 'Civilian '
$ws.Range("A7").Value = '     Civilian, New nominations'
$ws.Range("B7").Value = 223
$ws.Range("A8").Value = '     Civilian, Carryover nominations'
$ws.Range("B8").Value = 119
$ws.Range("A9").Value = '     Civilian, Confirmed  '
$ws.Range("B9").Value = 150
$ws.Range("A10").Value = '     Civilian, Unconfirmed  '
$ws.Range("B10").Value = 181
$ws.Range("A11").Value = '     Civilian, Withdrawn  '
$ws.Range("B11").Value = 11
$ws.Range("A12").Value = 'Civilian (FS, PHS, CG, NOAA)'
$ws.Range("A13").Value = '     Civilian (FS, PHS, CG, NOAA), New nominations'
$ws.Range("B13").Value = 1558
$ws.Range("A14").Value = '     Civilian (FS, PHS, CG, NOAA), Carryover nominations'
$ws.Range("B14").Value = 320
$ws.Range("A15").Value = '     Civilian (FS, PHS, CG, NOAA), Confirmed '
$ws.Range("B15").Value = 1335
$ws.Range("A16").Value = '     Civilian (FS, PHS, CG, NOAA), Unconfirmed '
$ws.Range("B16").Value = 543
$ws.Range("A17").Value = 'Air Force '
$ws.Range("A18").Value = '     Air Force, New nominations'
$ws.Range("B18").Value = 6213
$ws.Range("A19").Value = '     Air Force, Carryover nominations'
$ws.Range("B19").Value = 4952
$ws.Range("A20").Value = '     Air Force, Confirmed '
$ws.Range("B20").Value = 11159
$ws.Range("A21").Value = '     Air Force, Unconfirmed '
$ws.Range("B21").Value = 6
$ws.Range("A22").Value = 'Army '
$ws.Range("A23").Value = '     Army, New nominations'
$ws.Range("B23").Value = 8720
$ws.Range("A24").Value = '     Army, Carryover nominations'
$ws.Range("B24").Value = 2304
$ws.Range("A25").Value = '     Army, Confirmed '
$ws.Range("B25").Value = 11018
$ws.Range("A26").Value = '     Army, Unconfirmed '
$ws.Range("B26").Value = 6
$ws.Range("A27").Value = 'Navy '
$ws.Range("A28").Value = '     Navy, New nominations'
$ws.Range("B28").Value = 7165
$ws.Range("A29").Value = '     Navy, Carryover nominations'
$ws.Range("B29").Value = 21
$ws.Range("A30").Value = '     Navy, Confirmed '
$ws.Range("B30").Value = 7175
$ws.Range("A31").Value = '     Navy, Unconfirmed       '
$ws.Range("B31").Value = 11
$ws.Range("A32").Value = 'Marine Corps '
$ws.Range("A33").Value = '     Marine Corps, New nominations'
$ws.Range("B33").Value = 2332
$ws.Range("A34").Value = '     Marine Corps, Carryover nominations'
$ws.Range("B34").Value = 8
$ws.Range("A35").Value = '     Marine Corps, Confirmed '
$ws.Range("B35").Value = 2339
$ws.Range("A36").Value = '     Marine Corps, Unconfirmed  '
$ws.Range("A37").Value = 'Total new nominations'
# B37 was previously blank (the old "Summary" header row), so it has no
# number formatting yet -- borrow the "thousands separator, right
# aligned" look used by the other Total rows before putting the value in.
$ws.Range("B38").Copy()
$ws.Range("B37").PasteSpecial(-4122)
$ws.Range("B37").Value = 26211
$ws.Range("A38").Value = 'Total carryover nominations'
$ws.Range("B38").Value = 7724
$ws.Range("A39").Value = 'Total confirmed '
$ws.Range("B39").Value = 33176
$ws.Range("A40").Value = 'Total unconfirmed '
# B40 used to hold the "Total confirmed" number (thousands-separator
# style); the unconfirmed total below it instead uses the plain
# right-aligned style, so re-format it to match before writing the value.
$ws.Range("B41").Copy()
$ws.Range("B40").PasteSpecial(-4122)
$ws.Range("B40").Value = 748
$ws.Range("A41").Value = 'Total withdrawn '
$ws.Range("B41").Value = 11

# The old sheet had an extra trailing row (42) for "Total withdrawn";
# that label now lives on row 41, so remove the now-empty last row.
$ws.Rows.Item(42).Delete()
